$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits -----------------------------------------------------------
# E30: gun auto + burst time logged today (26+24+26 minutes)
$ws.Range("E30").Formula = "=(1/60)*(26+24+26)"

# New row 31 entry (trigger cd) - continues the shared "=(1/60)*(0)" pattern
# used by the rest of column B:E, plus the running daily total in F31.
$ws.Range("A31").NumberFormat = "[$-409]d\-mmm;@"
$ws.Range("A31").Value = 45318
$ws.Range("B31").Formula = "=(1/60)*(0)"
$ws.Range("C31").Formula = "=(1/60)*(0)"
$ws.Range("D31").Formula = "=(1/60)*(0)"
$ws.Range("E31").Formula = "=(1/60)*(0)"
$ws.Range("F31").Formula = "=SUM(B31:E31)"

# --- View / selection state -------------------------------------------
[void]$ws.Range("E31").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
